$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.86
$ws.Range("I2").Value = 4.33
$ws.Range("J2").Value = 2.6
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 8.5
$ws.Range("O2").Value = 1.33
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.63
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.75
$ws.Range("X2").Value = 8
$ws.Range("Z2").Value = 15
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 8.5
$ws.Range("AO2").Value = 10
$ws.Range("AR2").Value = 51
$ws.Range("AT2").Value = 2.63
$ws.Range("G3").Value = 1.86
$ws.Range("I3").Value = 5.25
$ws.Range("K3").Value = 1.8
$ws.Range("O3").Value = 1.63
$ws.Range("R3").Value = 1.33
$ws.Range("W3").Value = 4.5
$ws.Range("X3").Value = 7
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 21
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 126
$ws.Range("AI3").Value = 23
$ws.Range("AO3").Value = 11
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 34
$ws.Range("AY3").Value = 51
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 251
$ws.Range("G4").Value = 2.15
$ws.Range("K4").Value = 1.87
$ws.Range("M4").Value = 1.08
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.37
$ws.Range("R4").Value = 1.41
$ws.Range("G5").Value = 2.25
$ws.Range("K5").Value = 1.87
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.47
$ws.Range("P5").Value = 2.5
$ws.Range("Q5").Value = 2.6
$ws.Range("BD5").Value = 151
$ws.Range("G6").Value = 2.4
$ws.Range("H6").Value = 2.75
$ws.Range("K6").Value = 1.77
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 1.58
$ws.Range("R6").Value = 1.33
$ws.Range("H7").Value = 4.1
$ws.Range("I7").Value = 6
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 1.8
$ws.Range("X7").Value = 7
$ws.Range("AI7").Value = 34
$ws.Range("AM7").Value = 51
$ws.Range("AO7").Value = 7.5
$ws.Range("AU7").Value = 9
$ws.Range("AX7").Value = 34
$ws.Range("AY7").Value = 41
$ws.Range("AZ7").Value = 126
$ws.Range("BA7").Value = 151
$ws.Range("BB7").Value = 301
$ws.Range("G8").Value = 3.4
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 2.3
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 3.1
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 2.35
$ws.Range("R8").Value = 1.57
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("Y8").Value = 13
$ws.Range("AC8").Value = 7
$ws.Range("AE8").Value = 15
$ws.Range("AH8").Value = 6.5
$ws.Range("AJ8").Value = 10
$ws.Range("AP8").Value = 34
$ws.Range("AT8").Value = 2.38
$ws.Range("BA8").Value = 81
$ws.Range("BB8").Value = 251
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 5.5
